$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63: B63 - Pivot Turret text updated (Sherman -> turret wording, image renamed, button spacing)
$b63Text = @'
<Bold>e052 Pivot Turret</Bold> 
<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r8.24' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
If you want the turret to face a different sector, click tank counter on center of Battle Board. 
Alternatively, select buttons here:
<LineBreak/><LineBreak/>
                                   <InlineUIContainer><Button Content='  -  ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<InlineUIContainer><Image Name='c16TurretSherman75'  Height='150' Width='150'></Image></InlineUIContainer> 
<InlineUIContainer><Button Content='  +  ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
When you are satisfied with the current turret orientation, click turret image between buttons to continue.
'@
$ws.Range("B63").Value = $b63Text
$ws.Rows(63).RowHeight = 172.8

# Row 67: B67 - Rate of Fire text shortened
$b67Text = @'
<Bold>e053c Main Gun Firing - Rate of Fire</Bold> 
<InlineUIContainer><Button Content='r4.74.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r9.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<InlineUIContainer><Button Content='Rate of Fire' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
'@
$ws.Range("B67").Value = $b67Text
$ws.Rows(67).RowHeight = 72

# Row 69: B69 - MG Firing text reworded ("MG Gun Firing" -> "Machine Gun (MG) Firing")
$b69Text = @'
<Bold>e053e Machine Gun (MG) Firing - Select Target</Bold> 
<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r9.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
Click either the highlighted zone or an spotted target or <InlineUIContainer><Button Content='Skip MG' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> . 
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue53c' Height='100' Width='100'></Image></InlineUIContainer>
'@
$ws.Range("B69").Value = $b69Text

# Update active selection to match the saved view state (B68)
$ws.Range("B68").Select() | Out-Null
